$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: swap A1 and C1 (MLBSO00 <-> Date); B1 stays LNBSF00 ---
$ws.Range("A1").Value = "MLBSO00"
$ws.Range("C1").Value = "Date"

# --- Reset old data cells in columns A and B (rows 2-28) to blank / default style ---
$ws.Range("A2:B28").ClearContents()
$ws.Range("A2:B28").Style = "Normal"

# --- Move the date series into column C, keeping the date/time display style ---
$ws.Range("C2").Value = 45764
$ws.Range("C3").Value = 45763
$ws.Range("C4").Value = 45763
$ws.Range("C5").Value = 45762
$ws.Range("C6").Value = 45761
$ws.Range("C7").Value = 45758
$ws.Range("C8").Value = 45757
$ws.Range("C9").Value = 45756
$ws.Range("C10").Value = 45755
$ws.Range("C11").Value = 45754
$ws.Range("C12").Value = 45751
$ws.Range("C13").Value = 45750
$ws.Range("C14").Value = 45749
$ws.Range("C15").Value = 45748
$ws.Range("C16").Value = 45744
$ws.Range("C17").Value = 45743
$ws.Range("C18").Value = 45742
$ws.Range("C19").Value = 45741
$ws.Range("C20").Value = 45740
$ws.Range("C21").Value = 45737
$ws.Range("C22").Value = 45736
$ws.Range("C23").Value = 45735
$ws.Range("C24").Value = 45734
$ws.Range("C25").Value = 45733
$ws.Range("C26").Value = 45730
$ws.Range("C27").Value = 45729
$ws.Range("C28").Value = 45728
$ws.Range("C2:C28").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- New latest price observation lands in row 4, columns A and B ---
$ws.Range("A4").Value = 770.419
$ws.Range("B4").Value = 725.452
